# Generate Report for Handoff
# Adds a new row for "7b50f7ea-8d8c-4b4e-8a1b-36590b6108ef.md" to the
# Overview / zh-cn / de-de localization-status sheets, mirroring the
# existing "789ed9b5-a7c9-4a60-88d0-3392630c96f9.md" entry.

$wb = $excel.ActiveWorkbook

$fileName = "7b50f7ea-8d8c-4b4e-8a1b-36590b6108ef.md"
$pathName = "e2e\7b50f7ea-8d8c-4b4e-8a1b-36590b6108ef.md"
$ghUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7211b1676b027943e12f2945cb694b20ac28049f/e2e/7b50f7ea-8d8c-4b4e-8a1b-36590b6108ef.md"

# Leading "'" forces Excel to store the value as literal text instead of
# inferring a boolean/number/date/empty type (matches "True"/"False"/"" cells
# elsewhere in the sheet, which are plain text, not real booleans).
$emptyText = "'"
$trueText  = "'True"
$falseText = "'False"

# ---------------------------------------------------------------------
# Overview sheet: new row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = $fileName
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = $emptyText
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-18 18:41:57"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ghUrl, [Type]::Missing, [Type]::Missing, $pathName) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: new row 3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = $falseText
$wsZhCn.Cells.Item(3, 7).Value = "7b50f7ea-8d8c-4b4e-8a1b-36590b6108ef.e0b1a49999af74ce762d277a7d84a8efe402bd59.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value = "2016-08-18 18:41:52"
$wsZhCn.Cells.Item(3, 9).Value = $emptyText
$wsZhCn.Cells.Item(3, 10).Value = $emptyText
$wsZhCn.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 12).Value = $emptyText
$wsZhCn.Cells.Item(3, 13).Value = $trueText
$wsZhCn.Cells.Item(3, 14).Value = $emptyText
$wsZhCn.Cells.Item(3, 15).Value = $falseText
$wsZhCn.Cells.Item(3, 16).Value = $emptyText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $ghUrl, [Type]::Missing, [Type]::Missing, $fileName) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: new row 3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = $falseText
$wsDeDe.Cells.Item(3, 7).Value = "7b50f7ea-8d8c-4b4e-8a1b-36590b6108ef.e0b1a49999af74ce762d277a7d84a8efe402bd59.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value = "2016-08-18 18:41:57"
$wsDeDe.Cells.Item(3, 9).Value = $emptyText
$wsDeDe.Cells.Item(3, 10).Value = $emptyText
$wsDeDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 12).Value = $emptyText
$wsDeDe.Cells.Item(3, 13).Value = $trueText
$wsDeDe.Cells.Item(3, 14).Value = $emptyText
$wsDeDe.Cells.Item(3, 15).Value = $falseText
$wsDeDe.Cells.Item(3, 16).Value = $emptyText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $ghUrl, [Type]::Missing, [Type]::Missing, $fileName) | Out-Null

# ---------------------------------------------------------------------
# Resize the tables so the new rows are included (A1:P2 -> A1:P3 etc.)
# ---------------------------------------------------------------------
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))
$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P3"))
$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P3"))
